# Update the cryptos list worksheet with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new Price text, new Volume(1h) text). A $null entry
# means that column is left unchanged for that row.
$changes = @{
    2  = @("67.689.63", $null)
    3  = @("2.675.92", "  -1.20%  ")
    4  = @($null, "  +0.05%  ")
    5  = @("598.60", "  -0.21%  ")
    6  = @("166.39", "  +2.60%  ")
    7  = @($null, "  +0.02%  ")
    8  = @($null, "  +0.49%  ")
    9  = @("2.675.24", "  -1.18%  ")
    10 = @($null, "  +2.38%  ")
    11 = @($null, "  +1.13%  ")
    12 = @("0.362", "  +0.07%  ")
    13 = @($null, "  -1.65%  ")
    14 = @($null, "  -2.10%  ")
    15 = @("3.163.41", "  -0.96%  ")
    16 = @($null, "  -1.66%  ")
    17 = @("67.682.90", "  -1.46%  ")
    18 = @("2.675.51", "  -1.45%  ")
    19 = @($null, "  -1.05%  ")
    20 = @("7.74", "  +0.76%  ")
    21 = @("364.26", $null)
    22 = @("4.40", "  -3.32%  ")
    23 = @($null, "  -2.20%  ")
    24 = @($null, "  -3.90%  ")
    25 = @($null, "  +0.09%  ")
    26 = @("71.04", "  -4.24%  ")
    27 = @("10.12", "  +1.94%  ")
    29 = @($null, "  -2.80%  ")
    30 = @("1.00", "  -0.02%  ")
    31 = @("557.38", "  -7.08%  ")
    32 = @($null, "  -3.31%  ")
    33 = @($null, "  -3.51%  ")
    34 = @($null, "  -0.98%  ")
    35 = @($null, "  -1.73%  ")
    36 = @($null, "  +0.00%  ")
    37 = @($null, "  -4.63%  ")
    38 = @("19.54", "  -1.67%  ")
    39 = @("155.93", "  -2.76%  ")
    40 = @("0.373", "  -1.86%  ")
    41 = @("5.32", "  -2.24%  ")
    42 = @($null, "  -4.47%  ")
    43 = @($null, "  -0.38%  ")
    44 = @($null, "  -6.68%  ")
    45 = @($null, "  +0.02%  ")
    46 = @("40.33", "  -0.94%  ")
    47 = @($null, "  -5.65%  ")
    48 = @($null, "  -2.30%  ")
    49 = @("153.61", "  -2.99%  ")
    50 = @($null, "  -2.04%  ")
    51 = @($null, "  -3.02%  ")
}

# Price values that look like plain decimal numbers (e.g. "598.60") would be
# silently re-interpreted by Excel as a numeric value, losing the original
# text formatting (trailing zeros, fixed decimal places). Force those cells
# to Text format first so the literal string is preserved, matching the
# source data which stores prices as text.
function Test-LooksNumeric($s) {
    return $s -match '^[0-9]+(\.[0-9]+)?$'
}

foreach ($row in $changes.Keys) {
    $pair = $changes[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]

    if ($null -ne $newPrice) {
        $cell = $ws.Range("D$row")
        if (Test-LooksNumeric $newPrice) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $newPrice
    }
    if ($null -ne $newVolume) {
        $ws.Range("E$row").Value = $newVolume
    }
}
